$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.019.16"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "1.687.29"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "220.76"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "0.533"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "29.43"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").Value = "0.0640"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.929.49"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").Value = "1.671.40"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "10.25"
$ws.Range("E14").Value = "  +7.22%  "
$ws.Range("D15").Value = "0.608"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "4.15"
$ws.Range("E16").Value = "  +7.44%  "
$ws.Range("D17").Value = "31.049.74"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "66.90"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "247.80"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "10.04"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").Value = "158.64"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "15.95"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "3.52"
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").Value = "1.514.46"
$ws.Range("E34").Value = "  +6.31%  "
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "83.38"
$ws.Range("E37").Value = "  +10.10%  "
$ws.Range("D38").Value = "0.616"
$ws.Range("E38").Value = "  +10.04%  "
$ws.Range("E39").Value = "  +5.18%  "
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "0.0506"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "52.03"
$ws.Range("E47").Value = "  +7.03%  "
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").Value = "0.0₆0121"
$ws.Range("D51").Value = "93.85"
$ws.Range("E51").Value = "  +0.60%  "
